# STAGE 3 B1-B2 upload
# Populate the "B1" and "B2" evidence sheets with the new TxHash values.
# Shared-string entries must be created in this exact order so that the
# underlying sharedStrings table grows as: ..., B2!A2, B2!A3, B1!A2, B1!A3
# (matching the target workbook's uniqueCount 88 -> 92 append order).

$wb = $excel.ActiveWorkbook

$wsB1 = $wb.Worksheets.Item("B1")
$wsB2 = $wb.Worksheets.Item("B2")

# --- B2 gets its two new TxHash rows first (-> shared-string idx 88, 89) ---
$wsB2.Range("A2").Value = "63C72C2B34FF907482997050C7C87A1E65A2308C5C895B1D9DF9470A3674EDD0"
$wsB2.Range("A3").Value = "B3AE31C334892AAD882BE02E2DFBDE82F685EFC4425404D0A0F6350BD921D24A"

# --- B1 gets its two new TxHash rows next (-> shared-string idx 90, 91) ---
$wsB1.Range("A2").Value = "734218E2AC1A0F7E91FC306BF729ECEC36CC69819E18D12E81EA1D7F5CAD3FC5"
$wsB1.Range("A3").Value = "EA31E12D2B39908F31FD3CEE8C18951D6D2C1BD221D3493DA861AD7E7975EFDA"

# Best-effort column autosize so column A fits the new (longer) hash text.
# (Values chosen so the serialized column width lands as close as possible
# to the widths Excel's own best-fit produced: 85 for B1, ~82.86 for B2.)
$wsB1.Columns.Item(1).ColumnWidth = 84.16666666666667
$wsB2.Columns.Item(1).ColumnWidth = 82

# --- Update sheet selections / active tab ---
# B1 ends up with cell A4 selected (not the active tab).
$wsB1.Activate()
$wsB1.Range("A4").Select() | Out-Null

# B2 becomes the active/selected tab, with E30 selected, matching the
# workbook's new activeTab (index 22 = "B2").
$wsB2.Activate()
$wsB2.Range("E30").Select() | Out-Null
